$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 10, shifting existing rows 10..127 down to 11..128
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new record
$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(10, 3).Value = "Maule"
$ws.Cells.Item(10, 4).Value = Get-Date -Year 2022 -Month 3 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(10, 5).Value = 7
$ws.Cells.Item(10, 6).Value = 100112030
$ws.Cells.Item(10, 7).Value = "Poroto granado"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 400
$ws.Cells.Item(10, 11).Value = 23000
$ws.Cells.Item(10, 12).Value = 23000
$ws.Cells.Item(10, 13).Value = 23000
$ws.Cells.Item(10, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(10, 15).Value = "Región del Maule"
$ws.Cells.Item(10, 16).Value = 920
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
